# Apply edits described by the diff:
# - Q60: 2 -> 0
# - Q69: 2 -> 0
# - R1148, R1149: blank -> 0
# - Append new weekly rows 1150-1174 (data through 2024-12-16)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Two historical corrections in the "detect_structure" (Q) column.
$ws.Range("Q60").Value = 0
$ws.Range("Q69").Value = 0

# 2) The last two existing rows' "backup" (R) column now have a computed
#    value of 0 instead of being blank.
$ws.Range("R1148").Value = 0
$ws.Range("R1149").Value = 0

# 3) Newly appended weekly OHLCV rows. "Adj Close" (F) and "backup" (R)
#    are left blank for these rows, same as the source data.
$rows = @(
    @{ row=1150; A=45474; B=708; C=782.3499755859375; D=696.2999877929688; E=772.2000122070312; G=14122184; H=2024; I=7; J=1; K=0; L=0; M=0; N=27; O=0; P=0; Q=0 },
    @{ row=1151; A=45481; B=782.4000244140625; C=783.75; D=716.9000244140625; E=726.6500244140625; G=11796057; H=2024; I=7; J=8; K=0; L=0; M=0; N=28; O=1; P=0; Q=0 },
    @{ row=1152; A=45488; B=732; C=737.7999877929688; D=651.3499755859375; E=675.6500244140625; G=14597503; H=2024; I=7; J=15; K=0; L=0; M=0; N=29; O=0; P=0; Q=0 },
    @{ row=1153; A=45495; B=675; C=742; D=645; E=732.3499755859375; G=23462112; H=2024; I=7; J=22; K=0; L=0; M=0; N=30; O=2; P=0; Q=0 },
    @{ row=1154; A=45502; B=736.0999755859375; C=750.7000122070312; D=706.2000122070312; E=731.6500244140625; G=8130552; H=2024; I=7; J=29; K=0; L=0; M=0; N=31; O=0; P=0; Q=0 },
    @{ row=1155; A=45509; B=690.75; C=724; D=665.5499877929688; E=713.2000122070312; G=13910662; H=2024; I=8; J=5; K=0; L=0; M=0; N=32; O=0; P=0; Q=0 },
    @{ row=1156; A=45516; B=711; C=717.5; D=685.5; E=715.4500122070312; G=5537145; H=2024; I=8; J=12; K=0; L=0; M=0; N=33; O=0; P=0; Q=0 },
    @{ row=1157; A=45523; B=722; C=758.7999877929688; D=705.6500244140625; E=750.2000122070312; G=10589607; H=2024; I=8; J=19; K=0; L=0; M=0; N=34; O=0; P=0; Q=2 },
    @{ row=1158; A=45530; B=755.4500122070312; C=764.9000244140625; D=691.5; E=696.0999755859375; G=8733391; H=2024; I=8; J=26; K=0; L=0; M=0; N=35; O=0; P=0; Q=0 },
    @{ row=1159; A=45537; B=699; C=705; D=671.0499877929688; E=673.5499877929688; G=8155858; H=2024; I=9; J=2; K=0; L=0; M=0; N=36; O=0; P=0; Q=1 },
    @{ row=1160; A=45544; B=675.9000244140625; C=739; D=658.75; E=713.7000122070312; G=19700278; H=2024; I=9; J=9; K=0; L=0; M=0; N=37; O=2; P=0; Q=0 },
    @{ row=1161; A=45551; B=717.9500122070312; C=752; D=711.5499877929688; E=746.5; G=12615075; H=2024; I=9; J=16; K=0; L=0; M=0; N=38; O=0; P=0; Q=0 },
    @{ row=1162; A=45558; B=749.25; C=814.4000244140625; D=743.0999755859375; E=761.75; G=17364164; H=2024; I=9; J=23; K=0; L=0; M=0; N=39; O=0; P=0; Q=0 },
    @{ row=1163; A=45565; B=765.9000244140625; C=771.9500122070312; D=715.75; E=718.75; G=7596049; H=2024; I=9; J=30; K=0; L=0; M=0; N=40; O=0; P=0; Q=0 },
    @{ row=1164; A=45572; B=718; C=874.7000122070312; D=709.0499877929688; E=858.1500244140625; G=31275509; H=2024; I=10; J=7; K=0; L=0; M=0; N=41; O=1; P=0; Q=0 },
    @{ row=1165; A=45579; B=858.1500244140625; C=866.1500244140625; D=815; E=819; G=15095410; H=2024; I=10; J=14; K=0; L=0; M=0; N=42; O=0; P=0; Q=0 },
    @{ row=1166; A=45586; B=825; C=857.7000122070312; D=714.2999877929688; E=743.5; G=30406116; H=2024; I=10; J=21; K=0; L=0; M=0; N=43; O=0; P=0; Q=1 },
    @{ row=1167; A=45593; B=745.9500122070312; C=748; D=698.25; E=720.5999755859375; G=12069055; H=2024; I=10; J=28; K=0; L=0; M=0; N=44; O=0; P=0; Q=0 },
    @{ row=1168; A=45600; B=721; C=750; D=702.6500244140625; E=716.25; G=11369563; H=2024; I=11; J=4; K=0; L=0; M=0; N=45; O=0; P=0; Q=0 },
    @{ row=1169; A=45607; B=711; C=734.4000244140625; D=681.0999755859375; E=695.0999755859375; G=6581039; H=2024; I=11; J=11; K=0; L=0; M=0; N=46; O=2; P=0; Q=0 },
    @{ row=1170; A=45614; B=700; C=735; D=693.0499877929688; E=730.0499877929688; G=8138165; H=2024; I=11; J=18; K=0; L=0; M=0; N=47; O=0; P=0; Q=2 },
    @{ row=1171; A=45621; B=746.5999755859375; C=771; D=725.7000122070312; E=732.25; G=14955765; H=2024; I=11; J=25; K=0; L=0; M=0; N=48; O=0; P=0; Q=0 },
    @{ row=1172; A=45628; B=736.5999755859375; C=785.5; D=730.0499877929688; E=782.5; G=12069659; H=2024; I=12; J=2; K=0; L=0; M=0; N=49; O=0; P=0; Q=0 },
    @{ row=1173; A=45635; B=787.9500122070312; C=811.4000244140625; D=761; E=769; G=12730340; H=2024; I=12; J=9; K=0; L=0; M=0; N=50; O=0; P=0; Q=0 },
    @{ row=1174; A=45642; B=769.6500244140625; C=783; D=724.75; E=730.0499877929688; G=10422841; H=2024; I=12; J=16; K=0; L=0; M=0; N=51; O=0; P=0; Q=0 },
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.row, 1).Value = $r.A
    $ws.Cells.Item($r.row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r.row, 2).Value = $r.B
    $ws.Cells.Item($r.row, 3).Value = $r.C
    $ws.Cells.Item($r.row, 4).Value = $r.D
    $ws.Cells.Item($r.row, 5).Value = $r.E
    # Column F (Adj Close) intentionally left blank.
    $ws.Cells.Item($r.row, 7).Value = $r.G
    $ws.Cells.Item($r.row, 8).Value = $r.H
    $ws.Cells.Item($r.row, 9).Value = $r.I
    $ws.Cells.Item($r.row, 10).Value = $r.J
    $ws.Cells.Item($r.row, 11).Value = $r.K
    $ws.Cells.Item($r.row, 12).Value = $r.L
    $ws.Cells.Item($r.row, 13).Value = $r.M
    $ws.Cells.Item($r.row, 14).Value = $r.N
    $ws.Cells.Item($r.row, 15).Value = $r.O
    $ws.Cells.Item($r.row, 16).Value = $r.P
    $ws.Cells.Item($r.row, 17).Value = $r.Q
    # Column R (backup) intentionally left blank.
}

Write-Output "edit applied"
